# TimeSheet update: add 5 new rows (week of 18-22 May 2020) to the log table,
# extend Table1 to include the new rows, and update the selection.
#
# Shared-string table note: new unique text values are appended to
# xl/sharedStrings.xml in first-use order. The cells below are written in a
# specific order so the newly-created shared strings land at the same
# indices as in the target workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-DateCell($row) {
    $ws.Cells.Item(8, 1).Copy()
    $c = $ws.Cells.Item($row, 1)
    $c.PasteSpecial(-4122) | Out-Null   # xlPasteFormats - reuse existing date style
    return $c
}

function Set-TimeCell($row, $col) {
    $ws.Cells.Item(10, 2).Copy()
    $c = $ws.Cells.Item($row, $col)
    $c.PasteSpecial(-4122) | Out-Null   # xlPasteFormats - reuse existing time style
    return $c
}

# --- Write cells in the same order the original values were typed, so new
#     shared-string entries are appended in matching order. ------------------
(Set-TimeCell 16 2).Value = "00:00 - 03:00"
$ws.Cells.Item(16, 3).Value = "aanpassen van Analyse van de Opdracht V"

$ws.Cells.Item(19, 3).Value = "aanpassen analyse van de opdracht"
$ws.Cells.Item(19, 6).Value = "aanpassen van rescearch en oplossingen"

$ws.Cells.Item(20, 6).Value = "afwerken resaeach en oplossingen"

$ws.Cells.Item(17, 3).Value = "verbeteren desgne van oplossing"

(Set-TimeCell 18 2).Value = "16:00 - 17:00"
$ws.Cells.Item(18, 3).Value = "feedback design van aanpak6 analyse opdracht"

$ws.Cells.Item(20, 3).Value = "aanpassen research en ooplossingen"

$ws.Cells.Item(18, 5).Value = "sedric, stefan"

# --- Remaining cells (dates + cells reusing already-existing shared strings) -
(Set-DateCell 16).Value = "5/18/2020"
(Set-DateCell 17).Value = "5/19/2020"
(Set-DateCell 18).Value = "5/20/2020"
(Set-DateCell 19).Value = "5/21/2020"
(Set-DateCell 20).Value = "5/22/2020"

(Set-TimeCell 17 2).Value = "00:00 - 05:00"
(Set-TimeCell 19 2).Value = "00:00 - 03:00"
(Set-TimeCell 20 2).Value = "00:00 - 05:00"

$ws.Cells.Item(16, 5).Value = "Sedric Lodonou"
$ws.Cells.Item(17, 5).Value = "Sedric Lodonou"
$ws.Cells.Item(18, 4).Value = "Yvan"
$ws.Cells.Item(19, 5).Value = "Sedric Lodonou"
$ws.Cells.Item(20, 5).Value = "Sedric Lodonou"

$excel.CutCopyMode = $false

# --- Extend Table1 to cover the newly-added rows -------------------------
$table = $ws.ListObjects.Item("Table1")
$table.Resize($ws.Range("A5:G52"))

# --- Update selection to match the author's last-saved state -------------
$ws.Range("E18").Select()
